$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.531.35'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '1.640.81'
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.0000'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.64'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3789'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.84'
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3629'
$ws.Range("E9").Value = '  -1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08206'
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.236'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9978'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.54'
$ws.Range("E13").Value = '  -2.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.479'
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.400'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001243'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '1.633.59'
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.28'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06946'
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.610'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.56'
$ws.Range("E21").Value = '  -5.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.56'
$ws.Range("E23").Value = '  -3.56%  '
$ws.Range("D24").Value = '23.524.54'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.521'
$ws.Range("E25").Value = '  +3.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.078'
$ws.Range("E26").Value = '  -5.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.17'
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.96'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.274'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.49'
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("D31").Value = '1.816.55'
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.172'
$ws.Range("E32").Value = '  -7.18%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.632'
$ws.Range("E33").Value = '  -4.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.057'
$ws.Range("E34").Value = '  +8.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.45'
$ws.Range("E35").Value = '  +2.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02768'
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2500'
$ws.Range("E37").Value = '  -4.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08789'
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07142'
$ws.Range("E39").Value = '  -3.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.034'
$ws.Range("E40").Value = '  -6.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7147'
$ws.Range("E41").Value = '  -1.60%  '
$ws.Range("E42").Value = '  -3.01%  '
$ws.Range("E43").Value = '  -5.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.82'
$ws.Range("E44").Value = '  -4.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6584'
$ws.Range("E45").Value = '  -1.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9993'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.288'
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.980'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07986'
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '127.73'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.198'
$ws.Range("E51").Value = '  -2.92%  '
